# Horários.xlsx — add an extra "observações" column (C) with scheduling
# notes next to the existing availability column (B), and move the
# selection to C12.
#
# New shared strings must be introduced in this order so they land at the
# same shared-string-table indices as the authored workbook:
#   ... (19) "**"
#   (20) "o dia todo"
#   (21) "a partir das 18:45(pode variar de acordo com bus.. :) )"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the new column C (observações) to fit the longer note text.
$ws.Columns.Item(3).ColumnWidth = 47.1666666666666666

# Sabado 30/07 and Domingo 31/07 rows -> "o dia todo" (creates ss index 20)
$ws.Range("C8").Value = "o dia todo"
$ws.Range("C9").Value = "o dia todo"

# Quinta 28/07 and Sexta 29/07 rows -> the longer note (creates ss index 21)
$ws.Range("C6").Value = "a partir das 18:45(pode variar de acordo com bus.. :) )"
$ws.Range("C7").Value = "a partir das 18:45(pode variar de acordo com bus.. :) )"

# Move the active selection to C12, matching the saved view state.
$ws.Range("C12").Select() | Out-Null
